$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four account rows that were removed from the export, identified by the
# unique account number in column A (col 1). Using Find() for each one (rather
# than hard-coded row numbers) means each lookup is re-evaluated against the
# sheet's current state, so the rows can be removed in any order safely.
$accountsToDelete = @("005046919", "005063749", "004638738", "005219257")

foreach ($acct in $accountsToDelete) {
    $found = $ws.Columns.Item(1).Find($acct)
    if ($found -ne $null) {
        $found.EntireRow.Delete()
    }
}
